$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112095298
$ws.Range("B2").Value = 78713
$ws.Range("Q2").Value = 491104
$ws.Range("R2").Value = 6954282

# Row 3
$ws.Range("A3").Value = 112095473
$ws.Range("B3").Value = 78713
$ws.Range("Q3").Value = 491076
$ws.Range("R3").Value = 6954236

# Row 4
$ws.Range("A4").Value = 112095428
$ws.Range("B4").Value = 78713
$ws.Range("Q4").Value = 491096
$ws.Range("R4").Value = 6954259

# Row 5
$ws.Range("A5").Value = 112426713
$ws.Range("B5").Value = 78713
$ws.Range("Q5").Value = 490958
$ws.Range("R5").Value = 6953733
$ws.Range("Z5").Value = "18:39"
$ws.Range("AB5").Value = "18:39"

# Row 6
$ws.Range("A6").Value = 112426767
$ws.Range("B6").Value = 78713
$ws.Range("Q6").Value = 490949
$ws.Range("R6").Value = 6953753
$ws.Range("Z6").Value = "18:45"
$ws.Range("AB6").Value = "18:45"
